$wb = $excel.ActiveWorkbook

# Sheet2 starts empty (dimension A1:A1, no rows) -> populate A1:A10 with the
# numbers 0..9, growing the used range the same way the diff shows.
$ws2 = $wb.Worksheets.Item("Sheet2")
for ($i = 0; $i -lt 10; $i++) {
    $ws2.Cells.Item($i + 1, 1).Value = $i
}

# Sheet24 already holds 0..9 in A1:A10; re-apply the same values so the cells
# are (re)written on save, matching the diff's touch of those rows.
$ws24 = $wb.Worksheets.Item("Sheet24")
for ($i = 0; $i -lt 10; $i++) {
    $ws24.Cells.Item($i + 1, 1).Value = $i
}
